$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values on specific rows per repull of data
$ws.Range("F2").Value = -2
$ws.Range("F8").Value = 6
$ws.Range("F12").Value = -5
$ws.Range("F17").Value = -1
